# Fix km/h vs. m/s bug
#
# The "Geschwindigkeit" column on both sheets was labelled "(km/h)" but the
# pace (time) column formulas assumed the speed was already in m/s
# (800 or 1200 metres divided by speed divided by 86400 seconds/day).
# The fix: relabel the column "(m/s)" and convert the speed inputs from
# km/h to m/s by dividing by 3.6.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Both sheets are protected - unprotect so we can edit formulas/headers.
$ws1.Unprotect()
$ws2.Unprotect()

# --- Sheet 1: "6-7-8-9-10-11-12 kmh" -----------------------------------
# Relabel header (B9) from "Geschwindigkeit (km/h)" to "Geschwindigkeit (m/s)"
$ws1.Cells.Item(9, 2).Value = "Geschwindigkeit (m/s)"

# Convert the speed entry formulas from km/h to m/s (divide by 3.6)
$ws1.Cells.Item(10, 2).Formula = "=6/3.6"
$ws1.Cells.Item(11, 2).Formula = "=7/3.6"
$ws1.Cells.Item(12, 2).Formula = "=8/3.6"
$ws1.Cells.Item(13, 2).Formula = "=9/3.6"
$ws1.Cells.Item(14, 2).Formula = "=10/3.6"
$ws1.Cells.Item(15, 2).Formula = "=11/3.6"
$ws1.Cells.Item(16, 2).Formula = "=12/3.6"

# --- Sheet 2: "6.5-8-9.5-11-12.5-14 kmh" --------------------------------
$ws2.Cells.Item(9, 2).Value = "Geschwindigkeit (m/s)"

$ws2.Cells.Item(10, 2).Formula = "=6.5/3.6"
$ws2.Cells.Item(11, 2).Formula = "=8/3.6"
$ws2.Cells.Item(12, 2).Formula = "=9.5/3.6"
$ws2.Cells.Item(13, 2).Formula = "=11/3.6"
$ws2.Cells.Item(14, 2).Formula = "=12.5/3.6"
$ws2.Cells.Item(15, 2).Formula = "=14/3.6"

# --- Selection / active sheet bookkeeping (matches the authored edit) --
$ws1.Range("A11").Select()
$ws2.Activate()
$ws2.Range("B16").Select()
